$d = $word.ActiveDocument

# --- Insertion 1: new question before the "clumpy galaxies" question ---
$rng1 = $d.Content
$rng1.Find.Execute("Could you describe clumpy galaxies to us?") | Out-Null
$anchor1 = $rng1.Paragraphs(1)
$anchor1.Range.InsertParagraphBefore()
$anchor1.Range.Text = "Could you explain the KS and AD tests? What’s the difference? What does each test actually measure? (page 94)"

# --- Insertion 2: thirteen new questions after the "summarise this better in Chapter 3" question ---
$rng2 = $d.Content
$rng2.Find.Execute("summarise this better in Chapter 3") | Out-Null
$anchor2 = $rng2.Paragraphs(1)

$newQuestions = @(
  "Does the weighting scheme you have used in stellar mass and SFR keep your results valid? (page 101)",
  "How does balancing based on the stellar mass lead you to have equal number counts in each bin? (page 101)",
  "Could you summarise the error analysis you used from Cameron (2011)? What are you talking about with the beta function? (page 105)",
  "Fully break down and explain equation 3.1. You say it’s from Aird et al (2019), but you do not explain the fundamentals of this. (page 106)",
  "Are these results what are expected for the different stages of interaction? If so, why? (page 110)",
  "Why have you not used the confirmed merging galaxies in your pair sample? Could these not be assumed to be all at 0kpc separation? (page 110)",
  "You change your bin widths through the projected separation space. How does this affect your results? Is this something you can legitimately do? (page 110)",
  "How does taking the average measure the excess of the SFR from interaction? (page 110)",
  "What about only the pericenter and apocenter measurements on Figure 3.21? (page 111)",
  "How did you apply the Cameron (2011) methodology with respect to the projected separation? (page 112)",
  "Is the projected separation distribution with SFE for the full sample as expected? Why? (page 112)",
  "Should you not say something here about the timelines of the AGN and interaction? (page 113)",
  "Your weighting scheme in the AGN fractions is not clear… Did I actually do this two step fraction? (page 113)"
)

foreach ($q in $newQuestions) {
  $anchor2.Range.InsertParagraphAfter()
  $anchor2 = $anchor2.Next()
  $anchor2.Range.Text = $q
}

Write-Output "Paragraphs after edit:"
Write-Output $d.Paragraphs.Count
